$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column F previously held the PHP (OpenSwoole) comparison data; replace
# rows 2-11 with the new Java (Netty / Vert.x) comparison values.
$ws.Range("F2").Value  = "Netty / Vert.x"
$ws.Range("F3").Value  = "Event Loop + NIO"
$ws.Range("F4").Value  = "★★★★"
$ws.Range("F5").Value  = "~100k–500k"
$ws.Range("F6").Value  = "Trung bình–Cao"
$ws.Range("F7").Value  = "Trung bình (extension)"
$ws.Range("F8").Value  = "~20"
$ws.Range("F9").Value  = "~40–60"
$ws.Range("F10").Value = "Ổn định, enterprise, scale tốt"
$ws.Range("F11").Value = "JVM nặng, code verbose"

# Widen column F slightly to fit the new content (~40.16 characters).
$ws.Columns.Item(6).ColumnWidth = 39.330729166666664

# Move the active selection as it ended up after the edit.
$ws.Range("F16").Select()
